$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B -> C)
$ws.Columns.Item(2).Insert()

# Copy the header formatting from C1 (the original B1 header) onto the new B1 cell
$ws.Cells.Item(1, 3).Copy() | Out-Null
$ws.Cells.Item(1, 2).PasteSpecial(-4122) | Out-Null

# Set header row (row 1): A1, B1 (new), C1
$ws.Cells.Item(1, 1).Value = "Velocity_Bin"
$ws.Cells.Item(1, 2).Value = "Trening"
$ws.Cells.Item(1, 3).Value = "Acceleration_SMA"

# Data rows
$ws.Cells.Item(2, 1).Value = "10-15"
$ws.Cells.Item(2, 2).Value = "Duża Gra"
$ws.Cells.Item(2, 3).Value = 1.934777881311519

$ws.Cells.Item(3, 1).Value = "10-15"
$ws.Cells.Item(3, 2).Value = "Mała Gra"
$ws.Cells.Item(3, 3).Value = 3.27617809676885

$ws.Cells.Item(4, 1).Value = "5-10"
$ws.Cells.Item(4, 2).Value = "Duża Gra"
$ws.Cells.Item(4, 3).Value = 1.862376963063365

$ws.Cells.Item(5, 1).Value = "5-10"
$ws.Cells.Item(5, 2).Value = "Mała Gra"
$ws.Cells.Item(5, 3).Value = 2.847562606920276
